# Auto-generated Excel COM-interop edit script
# Implements: add rows 4-10 to "Logs" sheet, update "Dashboard" sheet rows 2-8,
# expand conditional formatting ranges, and update chart series source ranges.

$wb = $excel.ActiveWorkbook
$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# ---- 1. Append new rows 4-10 to the "Logs" sheet ----
# Row 4
$wsLogs.Range("A4").Value = "Sollicitatie marketingfunctie"
$wsLogs.Range("B4").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C4").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$wsLogs.Range("D4").Value = "Sollicitatie / Vacature"
$wsLogs.Range("E4").Value = "Geachte heer/mevrouw,`nHartelijk dank voor uw sollicitatie voor de functie van marketeer. Ik zal uw CV zorgvuldig bekijken en zal spoedig contact met u opnemen.`nMet vriendelijke groet,`n[Naam]"
$wsLogs.Range("F4").Value = "2025-06-22 17:33:14"
$wsLogs.Range("G4").Value = "Ja"

# Row 5
$wsLogs.Range("A5").Value = "Uitnodiging voor netwerkevent"
$wsLogs.Range("B5").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C5").Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$wsLogs.Range("D5").Value = "Uitnodiging / Evenement"
$wsLogs.Range("E5").Value = "Geachte heer/mevrouw,`nHartelijk dank voor de uitnodiging voor het zakelijke netwerkevent volgende maand. Ik waardeer uw uitnodiging en zal graag aanwezig zijn.`nMet vriendelijke groet,`n[Uw naam]"
$wsLogs.Range("F5").Value = "2025-06-22 17:33:17"
$wsLogs.Range("G5").Value = "Ja"

# Row 6
$wsLogs.Range("A6").Value = "Probleem met inloggen"
$wsLogs.Range("B6").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C6").Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$wsLogs.Range("D6").Value = "IT / Technisch probleem"
$wsLogs.Range("E6").Value = "Geachte klant,`nHartelijk dank voor uw bericht. Om u beter van dienst te kunnen zijn, kunnen we u vragen om enkele gegevens te verstrekken, zoals uw accountnaam en eventuele foutmeldingen die u krijgt. Met deze informatie kunnen we het probleem identificeren en hopelijk snel voor u oplossen. U kunt dit sturen naar support@email.com. `nMet vriendelijke groet,`n[Naam] `nKlantenservice medewerker"
$wsLogs.Range("F6").Value = "2025-06-22 17:33:19"
$wsLogs.Range("G6").Value = "Ja"

# Row 7
$wsLogs.Range("A7").Value = "Wat zijn jullie openingstijden?"
$wsLogs.Range("B7").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C7").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$wsLogs.Range("D7").Value = "Openingstijden / Locatie"
$wsLogs.Range("E7").Value = "Beste,`nBedankt voor je bericht. Onze openingstijden zijn van maandag tot en met vrijdag van 9.00 uur tot 18.00 uur. Op zaterdag zijn wij geopend van 10.00 uur tot 15.00 uur. Op zondag zijn wij gesloten.`nMet vriendelijke groet, `n[Naam Bedrijf]"
$wsLogs.Range("F7").Value = "2025-06-22 17:33:22"
$wsLogs.Range("G7").Value = "Ja"

# Row 8
$wsLogs.Range("A8").Value = "Probleem met inloggen"
$wsLogs.Range("B8").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C8").Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$wsLogs.Range("D8").Value = "IT / Technisch probleem"
$wsLogs.Range("E8").Value = "Beste klant,`nBedankt voor je bericht. Om je te kunnen helpen met het inlogprobleem op je account, hebben we wat meer informatie nodig. Kun je ons het e-mailadres en eventuele foutmeldingen die je krijgt sturen? We zullen ons best doen om zo snel mogelijk een oplossing te vinden.`nMet vriendelijke groet,`n[E-mailassistent]"
$wsLogs.Range("F8").Value = "2025-06-22 17:33:24"
$wsLogs.Range("G8").Value = "Ja"

# Row 9
$wsLogs.Range("A9").Value = "Offerte voor zakelijke samenwerking"
$wsLogs.Range("B9").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C9").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$wsLogs.Range("D9").Value = "Offerte / Prijsaanvraag"
$wsLogs.Range("F9").Value = "2025-06-22 17:33:25"
$wsLogs.Range("G9").Value = "Nee"

# Row 10
$wsLogs.Range("A10").Value = "Vragen over samenwerking"
$wsLogs.Range("B10").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C10").Value = "Kunnen we samenwerken aan een nieuw project?"
$wsLogs.Range("D10").Value = "Samenwerking / Partnerverzoek"
$wsLogs.Range("E10").Value = "Beste [Naam],`nDank voor je interesse in samenwerken aan een nieuw project. Graag ontvang ik meer details betreffende het project, zoals doelstellingen, tijdslijnen en verwachtingen. Op basis daarvan kunnen we een afspraak inplannen om de mogelijkheden te bespreken.`nMet vriendelijke groet,`n[Naam]"
$wsLogs.Range("F10").Value = "2025-06-22 17:33:27"
$wsLogs.Range("G10").Value = "Ja"

# ---- 2. Update the "Dashboard" sheet rows 2-8 (category summary) ----
# Row 2
$wsDash.Range("A2").Value = "IT / Technisch probleem"
$wsDash.Range("B2").Value = 3

# Row 3
$wsDash.Range("A3").Value = "Klacht / Probleem"
$wsDash.Range("B3").Value = 1

# Row 4
$wsDash.Range("A4").Value = "Sollicitatie / Vacature"
$wsDash.Range("B4").Value = 1

# Row 5
$wsDash.Range("A5").Value = "Uitnodiging / Evenement"
$wsDash.Range("B5").Value = 1

# Row 6
$wsDash.Range("A6").Value = "Openingstijden / Locatie"
$wsDash.Range("B6").Value = 1

# Row 7
$wsDash.Range("A7").Value = "Offerte / Prijsaanvraag"
$wsDash.Range("B7").Value = 1

# Row 8
$wsDash.Range("A8").Value = "Samenwerking / Partnerverzoek"
$wsDash.Range("B8").Value = 1

# ---- 3. Expand conditional formatting ranges on "Logs" sheet ----
# D2:D3 -> D2:D10  (Categorie column rules)
$fcsD = $wsLogs.Range("D2:D3").FormatConditions
for ($i = 1; $i -le $fcsD.Count(); $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($wsLogs.Range("D2:D10"))
}

# G2:G3 -> G2:G10  (Beantwoord column rules)
$fcsG = $wsLogs.Range("G2:G3").FormatConditions
for ($i = 1; $i -le $fcsG.Count(); $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($wsLogs.Range("G2:G10"))
}

# ---- 4. Update chart series source ranges (category totals now span A2:A8 / B2:B8) ----
$chartObj = $wsDash.ChartObjects().Item(1)
$chart = $chartObj.Chart()
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$8,Dashboard!`$B`$2:`$B`$8,1)"
